$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text (string) storage
# format instead of being auto-converted to numbers/percent by Excel.
$ws.Range("D2:E47").NumberFormat = "@"

$ws.Range("D2").Value = "274.66"
$ws.Range("E2").Value = "-1.31%"
$ws.Range("D3").Value = "26.84"
$ws.Range("E3").Value = "-2.11%"
$ws.Range("D4").Value = "4.909"
$ws.Range("E4").Value = "2.10%"
$ws.Range("D5").Value = "0.06328"
$ws.Range("E5").Value = "1.44%"
$ws.Range("D6").Value = "6.858"
$ws.Range("E6").Value = "-0.91%"
$ws.Range("D7").Value = "3.312"
$ws.Range("E7").Value = "1.33%"
$ws.Range("D8").Value = "1.254"
$ws.Range("E8").Value = "33.05%"
$ws.Range("D9").Value = "0.8695"
$ws.Range("E9").Value = "-1.20%"
$ws.Range("D10").Value = "0.1676"
$ws.Range("E10").Value = "15.40%"
$ws.Range("D11").Value = "0.05031"
$ws.Range("E11").Value = "-4.15%"
$ws.Range("D12").Value = "0.07486"
$ws.Range("E12").Value = "2.10%"
$ws.Range("D13").Value = "0.02968"
$ws.Range("E13").Value = "-4.37%"
$ws.Range("D14").Value = "0.09018"
$ws.Range("E14").Value = "-0.44%"
$ws.Range("D15").Value = "0.001571"
$ws.Range("E15").Value = "1.31%"
$ws.Range("D16").Value = "0.0006302"
$ws.Range("E16").Value = "0.57%"
$ws.Range("D17").Value = "0.005871"
$ws.Range("E17").Value = "-2.42%"
$ws.Range("D18").Value = "3.447"
$ws.Range("E18").Value = "-0.14%"
$ws.Range("D19").Value = "2.272"
$ws.Range("E19").Value = "-0.54%"
$ws.Range("D20").Value = "0.3138"
$ws.Range("E20").Value = "-0.26%"
$ws.Range("E21").Value = "2.49%"
$ws.Range("D22").Value = "3.906"
$ws.Range("E22").Value = "1.45%"
$ws.Range("D23").Value = "0.04350"
$ws.Range("E23").Value = "0.61%"
$ws.Range("D24").Value = "0.001175"
$ws.Range("E24").Value = "-0.21%"
$ws.Range("D25").Value = "0.004247"
$ws.Range("E25").Value = "-0.78%"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("D27").Value = "0.0001688"
$ws.Range("E27").Value = "-0.19%"
$ws.Range("D40").Value = "0.04044"
$ws.Range("E40").Value = "0.00%"
$ws.Range("D41").Value = "0.006728"
$ws.Range("E41").Value = "0.44%"
$ws.Range("D42").Value = "0.1165"
$ws.Range("E42").Value = "0.89%"
$ws.Range("D43").Value = "0.002192"
$ws.Range("E43").Value = "2.71%"
$ws.Range("D44").Value = "0.01073"
$ws.Range("E44").Value = "-11.46%"
$ws.Range("D45").Value = "0.00005303"
$ws.Range("E45").Value = "3.83%"
$ws.Range("D46").Value = "0.02102"
$ws.Range("E46").Value = "-29.63%"
$ws.Range("D47").Value = "1.490"
$ws.Range("E47").Value = "-37.29%"
